# Change age -> birthDate on the person class (students template)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Comment on A1 encodes the template's page range ":page\n<from>\n<to>";
# the end page moved from 7 to 9.
$comment = $ws.Range("A1").Comment
$null = $comment.Text(":page`n9`n5")

# Header label: "Age" -> "Birth Date"
$ws.Range("C2").Value = "Birth Date"

# Merge-field placeholder: {age} -> {birthDate|DATE}
$ws.Range("C3").Value = "{birthDate|DATE}"

# Active cell / selection moved from E7 to C10
[void]$ws.Range("C10").Select()
